# Swap the order of names in the "Recorded By" (column G) entries that
# currently read "System, dnasr281@gmail.com" so that they read
# "dnasr281@gmail.com, System" instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$target = "System, dnasr281@gmail.com"
$replacement = "dnasr281@gmail.com, System"

$searchRange = $ws.Range("G1:G319")

$first = $searchRange.Find($target)
if ($first -ne $null) {
    $firstAddress = $first.Address()
    $current = $first
    $continue = $true
    while ($continue) {
        $current.Value = $replacement
        $current = $searchRange.FindNext($current)
        if ($current -eq $null) {
            $continue = $false
        } elseif ($current.Address() -eq $firstAddress) {
            $continue = $false
        }
    }
}
